$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'28.573.50"
$ws.Range("D2").Style = $style
$style = $ws.Range("E2").Style
$ws.Range("E2").Value = "'  +2.03%  "
$ws.Range("E2").Style = $style
$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.864.11"
$ws.Range("D3").Style = $style
$style = $ws.Range("E3").Style
$ws.Range("E3").Value = "'  +1.94%  "
$ws.Range("E3").Style = $style
$style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = $style
$style = $ws.Range("E4").Style
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = $style
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'324.10"
$ws.Range("D5").Style = $style
$style = $ws.Range("E5").Style
$ws.Range("E5").Value = "'  -0.26%  "
$ws.Range("E5").Style = $style
$style = $ws.Range("E6").Style
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("E6").Style = $style
$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.4608"
$ws.Range("D7").Style = $style
$style = $ws.Range("E7").Style
$ws.Range("E7").Value = "'  -0.95%  "
$ws.Range("E7").Style = $style
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3872"
$ws.Range("D8").Style = $style
$style = $ws.Range("E8").Style
$ws.Range("E8").Value = "'  +0.19%  "
$ws.Range("E8").Style = $style
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.07864"
$ws.Range("D9").Style = $style
$style = $ws.Range("E9").Style
$ws.Range("E9").Value = "'  +0.11%  "
$ws.Range("E9").Style = $style
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.9728"
$ws.Range("D10").Style = $style
$style = $ws.Range("E10").Style
$ws.Range("E10").Value = "'  +1.50%  "
$ws.Range("E10").Style = $style
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'21.89"
$ws.Range("D11").Style = $style
$style = $ws.Range("E11").Style
$ws.Range("E11").Value = "'  +0.17%  "
$ws.Range("E11").Style = $style
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'1.870.88"
$ws.Range("D12").Style = $style
$style = $ws.Range("E12").Style
$ws.Range("E12").Value = "'  +1.92%  "
$ws.Range("E12").Style = $style
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'6.979"
$ws.Range("D13").Style = $style
$style = $ws.Range("E13").Style
$ws.Range("E13").Value = "'  +1.24%  "
$ws.Range("E13").Style = $style
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'5.684"
$ws.Range("D14").Style = $style
$style = $ws.Range("E14").Style
$ws.Range("E14").Value = "'  +0.17%  "
$ws.Range("E14").Style = $style
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.06931"
$ws.Range("D15").Style = $style
$style = $ws.Range("E15").Style
$ws.Range("E15").Value = "'  +1.00%  "
$ws.Range("E15").Style = $style
$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'88.22"
$ws.Range("D16").Style = $style
$style = $ws.Range("E16").Style
$ws.Range("E16").Value = "'  +1.23%  "
$ws.Range("E16").Style = $style
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = $style
$style = $ws.Range("E17").Style
$ws.Range("E17").Value = "'  +0.15%  "
$ws.Range("E17").Style = $style
$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'0.000009997"
$ws.Range("D18").Style = $style
$style = $ws.Range("E18").Style
$ws.Range("E18").Value = "'  +0.81%  "
$ws.Range("E18").Style = $style
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'16.79"
$ws.Range("D19").Style = $style
$style = $ws.Range("E19").Style
$ws.Range("E19").Value = "'  +1.19%  "
$ws.Range("E19").Style = $style
$style = $ws.Range("E20").Style
$ws.Range("E20").Value = "'  +0.16%  "
$ws.Range("E20").Style = $style
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'28.579.99"
$ws.Range("D21").Style = $style
$style = $ws.Range("E21").Style
$ws.Range("E21").Value = "'  +1.97%  "
$ws.Range("E21").Style = $style
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'5.262"
$ws.Range("D22").Style = $style
$style = $ws.Range("E22").Style
$ws.Range("E22").Value = "'  -0.96%  "
$ws.Range("E22").Style = $style
$style = $ws.Range("E23").Style
$ws.Range("E23").Value = "'  +0.81%  "
$ws.Range("E23").Style = $style
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'2.106"
$ws.Range("D24").Style = $style
$style = $ws.Range("E24").Style
$ws.Range("E24").Value = "'  +0.62%  "
$ws.Range("E24").Style = $style
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.091.75"
$ws.Range("D25").Style = $style
$style = $ws.Range("E25").Style
$ws.Range("E25").Value = "'  +1.46%  "
$ws.Range("E25").Style = $style
$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'152.54"
$ws.Range("D26").Style = $style
$style = $ws.Range("E26").Style
$ws.Range("E26").Value = "'  -0.77%  "
$ws.Range("E26").Style = $style
$style = $ws.Range("E27").Style
$ws.Range("E27").Value = "'  +0.61%  "
$ws.Range("E27").Style = $style
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'5.837"
$ws.Range("D28").Style = $style
$style = $ws.Range("E28").Style
$ws.Range("E28").Value = "'  +3.11%  "
$ws.Range("E28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'1.981"
$ws.Range("D29").Style = $style
$style = $ws.Range("E29").Style
$ws.Range("E29").Value = "'  +1.30%  "
$ws.Range("E29").Style = $style
$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'119.23"
$ws.Range("D30").Style = $style
$style = $ws.Range("E30").Style
$ws.Range("E30").Value = "'  +1.55%  "
$ws.Range("E30").Style = $style
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'0.09315"
$ws.Range("D31").Style = $style
$style = $ws.Range("E31").Style
$ws.Range("E31").Value = "'  +0.72%  "
$ws.Range("E31").Style = $style
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'0.9153"
$ws.Range("D32").Style = $style
$style = $ws.Range("E32").Style
$ws.Range("E32").Value = "'  -2.16%  "
$ws.Range("E32").Style = $style
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'5.269"
$ws.Range("D33").Style = $style
$style = $ws.Range("E33").Style
$ws.Range("E33").Value = "'  -0.06%  "
$ws.Range("E33").Style = $style
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'1.329"
$ws.Range("D34").Style = $style
$style = $ws.Range("E34").Style
$ws.Range("E34").Value = "'  +0.85%  "
$ws.Range("E34").Style = $style
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'3.320"
$ws.Range("D35").Style = $style
$style = $ws.Range("E35").Style
$ws.Range("E35").Value = "'  +0.84%  "
$ws.Range("E35").Style = $style
$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.05778"
$ws.Range("D36").Style = $style
$style = $ws.Range("E36").Style
$ws.Range("E36").Value = "'  -1.18%  "
$ws.Range("E36").Style = $style
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'1.147"
$ws.Range("D37").Style = $style
$style = $ws.Range("E37").Style
$ws.Range("E37").Value = "'  +0.75%  "
$ws.Range("E37").Style = $style
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.02063"
$ws.Range("D38").Style = $style
$style = $ws.Range("E38").Style
$ws.Range("E38").Value = "'  -2.74%  "
$ws.Range("E38").Style = $style
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'7.683"
$ws.Range("D39").Style = $style
$style = $ws.Range("E39").Style
$ws.Range("E39").Value = "'  -1.67%  "
$ws.Range("E39").Style = $style
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.5608"
$ws.Range("D40").Style = $style
$style = $ws.Range("E40").Style
$ws.Range("E40").Value = "'  +0.48%  "
$ws.Range("E40").Style = $style
$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.1779"
$ws.Range("D41").Style = $style
$style = $ws.Range("E41").Style
$ws.Range("E41").Value = "'  +1.12%  "
$ws.Range("E41").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'9.741"
$ws.Range("D42").Style = $style
$style = $ws.Range("E42").Style
$ws.Range("E42").Value = "'  -1.26%  "
$ws.Range("E42").Style = $style
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.07210"
$ws.Range("D43").Style = $style
$style = $ws.Range("E43").Style
$ws.Range("E43").Value = "'  +2.89%  "
$ws.Range("E43").Style = $style
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'11.71"
$ws.Range("D44").Style = $style
$style = $ws.Range("E44").Style
$ws.Range("E44").Value = "'  +0.41%  "
$ws.Range("E44").Style = $style
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.5278"
$ws.Range("D45").Style = $style
$style = $ws.Range("E45").Style
$ws.Range("E45").Value = "'  +0.40%  "
$ws.Range("E45").Style = $style
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.136"
$ws.Range("D46").Style = $style
$style = $ws.Range("E46").Style
$ws.Range("E46").Value = "'  +0.72%  "
$ws.Range("E46").Style = $style
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'1.127"
$ws.Range("D47").Style = $style
$style = $ws.Range("E47").Style
$ws.Range("E47").Value = "'  +1.60%  "
$ws.Range("E47").Style = $style
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'1.832"
$ws.Range("D48").Style = $style
$style = $ws.Range("E48").Style
$ws.Range("E48").Value = "'  +0.28%  "
$ws.Range("E48").Style = $style
$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'112.62"
$ws.Range("D49").Style = $style
$style = $ws.Range("E49").Style
$ws.Range("E49").Value = "'  -0.03%  "
$ws.Range("E49").Style = $style
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'2.407"
$ws.Range("D50").Style = $style
$style = $ws.Range("E50").Style
$ws.Range("E50").Value = "'  +3.66%  "
$ws.Range("E50").Style = $style
$style = $ws.Range("E51").Style
$ws.Range("E51").Value = "'  +0.21%  "
$ws.Range("E51").Style = $style
